# Update Il27-Il27ra.xlsx with new TPM-derived values.
#
# Net effect vs. the original workbook:
#  - The data block shrinks from 6 rows (ECs-as-sender x3, FAPs-as-sender x3)
#    down to 3 rows (FAPs-as-sender only), so the used range becomes A1:T4.
#  - The numeric payload for the surviving FAPs rows is refreshed with new
#    TPM-derived figures.
#  - Because the "ECs" sending-cluster rows are dropped and the cluster-name
#    strings are rewritten, the shared-string table ends up re-interning
#    FAPs/Il27/Il27ra/ECs/MuSCs in a new order (FAPs, Il27, Il27ra, ECs,
#    MuSCs) - we reproduce that by first blanking every cell that names one
#    of those clusters (so the old strings fall out of the shared-string
#    pool entirely), then writing them back in the desired order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Blank every cell in the A:D / rows 2-7 block so none of the cluster-name
#    strings (ECs, FAPs, Il27, Il27ra, MuSCs) are referenced anywhere left in
#    the workbook.
$ws.Range("A2:D7").ClearContents()

# 2) Re-populate rows 2-4 (sending cluster / ligand / receptor / target
#    cluster columns) in the order that makes the shared strings intern as
#    FAPs, Il27, Il27ra, ECs, MuSCs.
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Il27"
$ws.Range("C2").Value = "Il27ra"
$ws.Range("D2").Value = "ECs"

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Il27"
$ws.Range("C3").Value = "Il27ra"
$ws.Range("D3").Value = "FAPs"

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Il27"
$ws.Range("C4").Value = "Il27ra"
$ws.Range("D4").Value = "MuSCs"

# 3) Refresh the numeric payload (columns E-T) for the 3 surviving rows.
$data = @(
    @(1,0.3333333333333333,0.1127286666666667,0.338186,1,1,3,1,0.5127123333333333,1.538137,0.4737471074894348,0.4737471074894348,0.05779737772022222,0.5201763994819999,0.4737471074894348,0.4737471074894348),
    @(1,0.3333333333333333,0.1127286666666667,0.338186,1,1,3,1,0.5455243333333334,1.636573,0.5040654538219332,0.5040654538219331,0.06149623073088889,0.553466076578,0.5040654538219332,0.5040654538219331),
    @(1,0.3333333333333333,0.1127286666666667,0.338186,1,1,2,0.6666666666666666,0.02401233333333333,0.07203699999999999,0.02218743868863203,0.02218743868863203,0.002706878320222222,0.02436190488199999,0.02218743868863203,0.02218743868863203)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        # Columns E..T start at column index 5.
        $ws.Cells.Item($rowNum, $c + 5).Value = $rowVals[$c]
    }
}

# 4) Remove the now-obsolete rows 5-7, shrinking the used range to A1:T4.
$ws.Range("A5:T7").EntireRow.Delete()
